$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.02"
$ws.Range("E2").Value = "'-0.44%"
$ws.Range("D3").Value = "'38.85"
$ws.Range("E3").Value = "'7.05%"
$ws.Range("D4").Value = "'5.108"
$ws.Range("E4").Value = "'1.01%"
$ws.Range("D5").Value = "'0.08082"
$ws.Range("E5").Value = "'-0.50%"
$ws.Range("D6").Value = "'1.933"
$ws.Range("E6").Value = "'-2.49%"
$ws.Range("D7").Value = "'4.182"
$ws.Range("E7").Value = "'0.48%"
$ws.Range("D8").Value = "'7.995"
$ws.Range("D9").Value = "'0.9306"
$ws.Range("E9").Value = "'0.16%"
$ws.Range("D10").Value = "'0.1461"
$ws.Range("E10").Value = "'-0.20%"
$ws.Range("D11").Value = "'0.1921"
$ws.Range("E11").Value = "'-0.10%"
$ws.Range("D12").Value = "'0.09066"
$ws.Range("E12").Value = "'-0.88%"
$ws.Range("D13").Value = "'0.03511"
$ws.Range("E13").Value = "'1.99%"
$ws.Range("D14").Value = "'0.09787"
$ws.Range("E14").Value = "'-1.00%"
$ws.Range("D15").Value = "'0.001391"
$ws.Range("E15").Value = "'-2.41%"
$ws.Range("D16").Value = "'0.005844"
$ws.Range("E16").Value = "'-13.77%"
$ws.Range("D17").Value = "'3.776"
$ws.Range("E17").Value = "'-1.44%"
$ws.Range("E18").Value = "'0.54%"
$ws.Range("D19").Value = "'0.3443"
$ws.Range("E19").Value = "'-0.38%"
$ws.Range("E20").Value = "'5.06%"
$ws.Range("D21").Value = "'4.680"
$ws.Range("E21").Value = "'-3.18%"
$ws.Range("E22").Value = "'3.17%"
$ws.Range("D23").Value = "'0.04371"
$ws.Range("E23").Value = "'-0.28%"
$ws.Range("E24").Value = "'0.20%"
$ws.Range("D25").Value = "'0.004274"
$ws.Range("E25").Value = "'2.30%"
$ws.Range("D26").Value = "'0.0001303"
$ws.Range("E26").Value = "'-0.05%"
$ws.Range("D39").Value = "'0.02035"
$ws.Range("E39").Value = "'-0.08%"
$ws.Range("D40").Value = "'0.05059"
$ws.Range("E40").Value = "'-1.47%"
$ws.Range("D41").Value = "'0.007529"
$ws.Range("E41").Value = "'0.63%"
$ws.Range("D42").Value = "'0.009742"
$ws.Range("E42").Value = "'-4.01%"
$ws.Range("E43").Value = "'-2.07%"
$ws.Range("D44").Value = "'0.002126"
$ws.Range("E44").Value = "'1.38%"
$ws.Range("D45").Value = "'0.009925"
$ws.Range("E45").Value = "'0.37%"
$ws.Range("D46").Value = "'0.00006199"
$ws.Range("E46").Value = "'-1.76%"
$ws.Range("E47").Value = "'-0.03%"
$ws.Range("D48").Value = "'0.002876"
$ws.Range("D49").Value = "'0.001803"
$ws.Range("E49").Value = "'12.38%"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'-0.03%"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'-0.03%"
